$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting existing C..I to D..J
$ws.Columns("C:C").Insert()

# New column C header + first data row
$ws.Range("C1").Value = "00 (utility)"
$ws.Range("C2").Value = "abyss_tile"

# New legend entry in column A
$ws.Range("A6").Value = "stained_glass"

# New item in (formerly-C, now-D) "01 (traps)" column
$ws.Range("D7").Value = "rusty_platform_double"

# New pink highlight style (fill FFFF66CC), matching the style used elsewhere
# for colour-coded legend / category cells
$ws.Range("C2").Interior.Color = 13395711
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").HorizontalAlignment = -4108

$ws.Range("A6").Interior.Color = 13395711
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").HorizontalAlignment = -4108

# D7 reuses the existing "material" style (same as A2 / D2)
$ws.Range("D7").Interior.Color = $ws.Range("D2").Interior.Color
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").HorizontalAlignment = -4108

# Column width / sheet layout tweaks
$ws.Range("B1").ColumnWidth = 7.15
$ws.Range("C1").ColumnWidth = 21.15

# Selection moved by the author
$ws.Range("D14").Select()
